$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.72976833333333
$ws.Range("H2").Value = 92.18930499999999
$ws.Range("I2").Value = 0.5616793902924558
$ws.Range("J2").Value = 0.5616793902924558
$ws.Range("M2").Value = 30.72976833333333
$ws.Range("N2").Value = 92.18930499999999
$ws.Range("O2").Value = 0.5616793902924558
$ws.Range("P2").Value = 0.5616793902924558
$ws.Range("Q2").Value = 944.3186618203358
$ws.Range("R2").Value = 8498.867956383023
$ws.Range("S2").Value = 0.3154837374793049
$ws.Range("T2").Value = 0.3154837374793049

$ws.Range("G3").Value = 30.72976833333333
$ws.Range("H3").Value = 92.18930499999999
$ws.Range("I3").Value = 0.5616793902924558
$ws.Range("J3").Value = 0.5616793902924558
$ws.Range("M3").Value = 4.690023666666666
$ws.Range("O3").Value = 0.08572435707863904
$ws.Range("P3").Value = 0.08572435707863904
$ws.Range("Q3").Value = 144.1233407545172
$ws.Range("R3").Value = 1297.110066790655
$ws.Range("S3").Value = 0.04814960461714275
$ws.Range("T3").Value = 0.04814960461714275

$ws.Range("G4").Value = 30.72976833333333
$ws.Range("H4").Value = 92.18930499999999
$ws.Range("I4").Value = 0.5616793902924558
$ws.Range("J4").Value = 0.5616793902924558
$ws.Range("M4").Value = 19.25278733333333
$ws.Range("N4").Value = 57.758362
$ws.Range("O4").Value = 0.3519028758536682
$ws.Range("P4").Value = 0.3519028758536682
$ws.Range("Q4").Value = 591.6336945242678
$ws.Range("R4").Value = 5324.703250718409
$ws.Range("S4").Value = 0.1976565927516501
$ws.Range("T4").Value = 0.1976565927516501

$ws.Range("G5").Value = 30.72976833333333
$ws.Range("H5").Value = 92.18930499999999
$ws.Range("I5").Value = 0.5616793902924558
$ws.Range("J5").Value = 0.5616793902924558
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.037935
$ws.Range("N5").Value = 0.113805
$ws.Range("O5").Value = 0.0006933767752369208
$ws.Range("P5").Value = 0.0006933767752369208
$ws.Range("Q5").Value = 1.165733761725
$ws.Range("R5").Value = 10.491603855525
$ws.Range("S5").Value = 0.0003894554443580228
$ws.Range("T5").Value = 0.0003894554443580228

$ws.Range("G6").Value = 4.690023666666666
$ws.Range("I6").Value = 0.08572435707863904
$ws.Range("J6").Value = 0.08572435707863904
$ws.Range("M6").Value = 30.72976833333333
$ws.Range("N6").Value = 92.18930499999999
$ws.Range("O6").Value = 0.5616793902924558
$ws.Range("P6").Value = 0.5616793902924558
$ws.Range("Q6").Value = 144.1233407545172
$ws.Range("R6").Value = 1297.110066790655
$ws.Range("S6").Value = 0.04814960461714275
$ws.Range("T6").Value = 0.04814960461714275

$ws.Range("G7").Value = 4.690023666666666
$ws.Range("I7").Value = 0.08572435707863904
$ws.Range("J7").Value = 0.08572435707863904
$ws.Range("M7").Value = 4.690023666666666
$ws.Range("O7").Value = 0.08572435707863904
$ws.Range("P7").Value = 0.08572435707863904
$ws.Range("Q7").Value = 21.99632199389344
$ws.Range("R7").Value = 197.966897945041
$ws.Range("S7").Value = 0.007348665396546011
$ws.Range("T7").Value = 0.007348665396546011

$ws.Range("G8").Value = 4.690023666666666
$ws.Range("I8").Value = 0.08572435707863904
$ws.Range("J8").Value = 0.08572435707863904
$ws.Range("M8").Value = 19.25278733333333
$ws.Range("N8").Value = 57.758362
$ws.Range("O8").Value = 0.3519028758536682
$ws.Range("P8").Value = 0.3519028758536682
$ws.Range("Q8").Value = 90.29602824263354
$ws.Range("R8").Value = 812.6642541837019
$ws.Range("S8").Value = 0.03016664778667984
$ws.Range("T8").Value = 0.03016664778667984

$ws.Range("G9").Value = 4.690023666666666
$ws.Range("I9").Value = 0.08572435707863904
$ws.Range("J9").Value = 0.08572435707863904
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.037935
$ws.Range("N9").Value = 0.113805
$ws.Range("O9").Value = 0.0006933767752369208
$ws.Range("P9").Value = 0.0006933767752369208
$ws.Range("Q9").Value = 0.1779160477949999
$ws.Range("R9").Value = 1.601244430155
$ws.Range("S9").Value = 0.00005943927827044504
$ws.Range("T9").Value = 0.00005943927827044504

$ws.Range("G10").Value = 19.25278733333333
$ws.Range("H10").Value = 57.758362
$ws.Range("I10").Value = 0.3519028758536682
$ws.Range("J10").Value = 0.3519028758536682
$ws.Range("M10").Value = 30.72976833333333
$ws.Range("N10").Value = 92.18930499999999
$ws.Range("O10").Value = 0.5616793902924558
$ws.Range("P10").Value = 0.5616793902924558
$ws.Range("Q10").Value = 591.6336945242678
$ws.Range("R10").Value = 5324.703250718409
$ws.Range("S10").Value = 0.1976565927516501
$ws.Range("T10").Value = 0.1976565927516501

$ws.Range("G11").Value = 19.25278733333333
$ws.Range("H11").Value = 57.758362
$ws.Range("I11").Value = 0.3519028758536682
$ws.Range("J11").Value = 0.3519028758536682
$ws.Range("M11").Value = 4.690023666666666
$ws.Range("O11").Value = 0.08572435707863904
$ws.Range("P11").Value = 0.08572435707863904
$ws.Range("Q11").Value = 90.29602824263354
$ws.Range("R11").Value = 812.6642541837019
$ws.Range("S11").Value = 0.03016664778667984
$ws.Range("T11").Value = 0.03016664778667984

$ws.Range("G12").Value = 19.25278733333333
$ws.Range("H12").Value = 57.758362
$ws.Range("I12").Value = 0.3519028758536682
$ws.Range("J12").Value = 0.3519028758536682
$ws.Range("M12").Value = 19.25278733333333
$ws.Range("N12").Value = 57.758362
$ws.Range("O12").Value = 0.3519028758536682
$ws.Range("P12").Value = 0.3519028758536682
$ws.Range("Q12").Value = 370.6698201025605
$ws.Range("R12").Value = 3336.028380923044
$ws.Range("S12").Value = 0.1238356340340822
$ws.Range("T12").Value = 0.1238356340340822

$ws.Range("G13").Value = 19.25278733333333
$ws.Range("H13").Value = 57.758362
$ws.Range("I13").Value = 0.3519028758536682
$ws.Range("J13").Value = 0.3519028758536682
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.037935
$ws.Range("N13").Value = 0.113805
$ws.Range("O13").Value = 0.0006933767752369208
$ws.Range("P13").Value = 0.0006933767752369208
$ws.Range("Q13").Value = 0.73035448749
$ws.Range("R13").Value = 6.57319038741
$ws.Range("S13").Value = 0.0002440012812560149
$ws.Range("T13").Value = 0.0002440012812560149

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.037935
$ws.Range("H14").Value = 0.113805
$ws.Range("I14").Value = 0.0006933767752369208
$ws.Range("J14").Value = 0.0006933767752369208
$ws.Range("M14").Value = 30.72976833333333
$ws.Range("N14").Value = 92.18930499999999
$ws.Range("O14").Value = 0.5616793902924558
$ws.Range("P14").Value = 0.5616793902924558
$ws.Range("Q14").Value = 1.165733761725
$ws.Range("R14").Value = 10.491603855525
$ws.Range("S14").Value = 0.0003894554443580228
$ws.Range("T14").Value = 0.0003894554443580228

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.037935
$ws.Range("H15").Value = 0.113805
$ws.Range("I15").Value = 0.0006933767752369208
$ws.Range("J15").Value = 0.0006933767752369208
$ws.Range("M15").Value = 4.690023666666666
$ws.Range("O15").Value = 0.08572435707863904
$ws.Range("P15").Value = 0.08572435707863904
$ws.Range("Q15").Value = 0.1779160477949999
$ws.Range("R15").Value = 1.601244430155
$ws.Range("S15").Value = 0.00005943927827044504
$ws.Range("T15").Value = 0.00005943927827044504

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.037935
$ws.Range("H16").Value = 0.113805
$ws.Range("I16").Value = 0.0006933767752369208
$ws.Range("J16").Value = 0.0006933767752369208
$ws.Range("M16").Value = 19.25278733333333
$ws.Range("N16").Value = 57.758362
$ws.Range("O16").Value = 0.3519028758536682
$ws.Range("P16").Value = 0.3519028758536682
$ws.Range("Q16").Value = 0.73035448749
$ws.Range("R16").Value = 6.57319038741
$ws.Range("S16").Value = 0.0002440012812560149
$ws.Range("T16").Value = 0.0002440012812560149

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.037935
$ws.Range("H17").Value = 0.113805
$ws.Range("I17").Value = 0.0006933767752369208
$ws.Range("J17").Value = 0.0006933767752369208
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.037935
$ws.Range("N17").Value = 0.113805
$ws.Range("O17").Value = 0.0006933767752369208
$ws.Range("P17").Value = 0.0006933767752369208
$ws.Range("Q17").Value = 0.001439064225
$ws.Range("R17").Value = 0.012951578025
$ws.Range("S17").Value = 0.0000004807713524379514
$ws.Range("T17").Value = 0.0000004807713524379514
